# Fixed TBD issue and regex issue
# Replace the placeholder "nan" used for unassigned instructors/rooms with "TBD".
# Only text that is part of a course-schedule entry (e.g. "EG-208-03804-nan-IDE-214")
# should change; the handful of cells whose entire value is just the literal
# string "nan" (used elsewhere as a generic "not available" marker) must stay
# as-is. Every schedule-entry occurrence of the placeholder is immediately
# preceded by a hyphen ("-nan"), while the standalone cells are exactly "nan"
# with nothing before it, so replacing "-nan" -> "-TBD" everywhere is safe and
# precise.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
foreach ($cell in $used.Cells) {
    $val = $cell.Value2
    if ($null -ne $val -and $val -is [string] -and $val.Contains("-nan")) {
        $cell.Value = $val.Replace("-nan", "-TBD")
    }
}
